# Updated symbol list on Mon Feb 13 20:27:53 UTC 2023 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures for the existing coin rows,
# and swaps the BOLO / CoinbaseStockToken rows (48 <-> 49) with their new data.
#
# All Price/Volume values are stored as text (matching the source sheet's
# inlineStr cells), so numeric-looking strings are entered with a leading
# apostrophe to force text storage instead of Excel auto-converting them
# to numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'289.04"
$ws.Range("E2").Value = "'-9.54%"
$ws.Range("D3").Value = "'40.46"
$ws.Range("E3").Value = "'-2.49%"
$ws.Range("D4").Value = "'5.042"
$ws.Range("E4").Value = "'-4.09%"
$ws.Range("E5").Value = "'-5.85%"
$ws.Range("D6").Value = "'4.281"
$ws.Range("E6").Value = "'-1.71%"
$ws.Range("D7").Value = "'1.547"
$ws.Range("E7").Value = "'-11.17%"
$ws.Range("D8").Value = "'0.9183"
$ws.Range("E8").Value = "'-2.85%"
$ws.Range("D9").Value = "'0.1168"
$ws.Range("E9").Value = "'-5.74%"
$ws.Range("D10").Value = "'0.1734"
$ws.Range("E10").Value = "'-7.07%"
$ws.Range("D11").Value = "'0.08688"
$ws.Range("E11").Value = "'-6.03%"
$ws.Range("D12").Value = "'0.04164"
$ws.Range("E12").Value = "'1.33%"
$ws.Range("D13").Value = "'0.1051"
$ws.Range("E13").Value = "'-0.03%"
$ws.Range("D14").Value = "'0.001272"
$ws.Range("E14").Value = "'-1.31%"
$ws.Range("D15").Value = "'0.005808"
$ws.Range("E15").Value = "'0.52%"
$ws.Range("D16").Value = "'3.394"
$ws.Range("E16").Value = "'1.23%"
$ws.Range("E18").Value = "'-1.20%"
$ws.Range("D19").Value = "'7.835"
$ws.Range("E19").Value = "'-6.67%"
$ws.Range("D20").Value = "'0.1351"
$ws.Range("E20").Value = "'-0.13%"
$ws.Range("D21").Value = "'0.2884"
$ws.Range("E21").Value = "'2.18%"
$ws.Range("D22").Value = "'0.03872"
$ws.Range("E22").Value = "'-3.76%"
$ws.Range("D23").Value = "'0.001268"
$ws.Range("E23").Value = "'-0.13%"
$ws.Range("D24").Value = "'0.003851"
$ws.Range("E24").Value = "'-6.55%"
$ws.Range("D25").Value = "'0.0001282"
$ws.Range("E25").Value = "'0.77%"
$ws.Range("D26").Value = "'0.0003725"
$ws.Range("E26").Value = "'-95.03%"
$ws.Range("D38").Value = "'0.02315"
$ws.Range("E38").Value = "'-9.74%"
$ws.Range("D39").Value = "'0.04968"
$ws.Range("E39").Value = "'-7.03%"
$ws.Range("D40").Value = "'0.006728"
$ws.Range("E40").Value = "'237.91%"
$ws.Range("D41").Value = "'0.007683"
$ws.Range("E41").Value = "'-1.38%"
$ws.Range("D42").Value = "'0.1273"
$ws.Range("E42").Value = "'-3.41%"
$ws.Range("D43").Value = "'0.007380"
$ws.Range("E43").Value = "'5.04%"
$ws.Range("D44").Value = "'0.007067"
$ws.Range("E44").Value = "'-14.35%"
$ws.Range("D45").Value = "'0.3122"
$ws.Range("E45").Value = "'-1.49%"
$ws.Range("D46").Value = "'0.00006442"
$ws.Range("E46").Value = "'-3.79%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.004202"
$ws.Range("E48").Value = "'-0.01%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.02982"
$ws.Range("E49").Value = "'-85.06%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.05%"
